# Applies the "products.xlsx" catalog edit described by the commit:
#   "FIxed real data, changed images to load with long names :)"
#
# Summary of the change:
#  1) Row 95 (product #94, "Dschenresi ... mit Mantra 30cm") is moved to the
#     bottom of the sheet (new row 240) and its image-filename list is
#     corrected (the bogus "MBAva018.JPG" entry is dropped).
#  2) Row 148's image filenames get a casing fix
#     ("TBT020.jpg"/"TBT018.jpg" -> "TBT020.JPG"/"TBT018.JPG").
#  3) Row 178 (product #177, "Gruene Tara,21 Taras") turns out to not be a
#     real/found product, so it is moved out to a brand new worksheet named
#     "catalog_not_found".
#  4) Small view/selection bookkeeping (active cell, which tab is selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Move row 95 to row 240 (end of sheet) and fix its image list text.
# ---------------------------------------------------------------------
$ws.Range("A95:G95").Cut($ws.Range("A240:G240"))
$ws.Cells.Item(240, 6).Value = "MBAva014.JPG, MBAva015.JPG, MBAva016.JPG, MBAva017.JPG, MBAva019.JPG"

# ---------------------------------------------------------------------
# 2) Fix casing of the image list referenced by row 148.
# ---------------------------------------------------------------------
$ws.Cells.Item(148, 6).Value = "TBT020.JPG, TBT019.JPG, TBT018.JPG"

# ---------------------------------------------------------------------
# 3) Add the "catalog_not_found" worksheet after the existing sheet and
#    move row 178 into it as its first (and only) data row.
# ---------------------------------------------------------------------
$notFound = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$notFound.Name = "catalog_not_found"

$ws.Range("A178:G178").Cut($notFound.Range("A1:G1"))

# ---------------------------------------------------------------------
# 4) View bookkeeping: leave the original sheet active/selected, with the
#    selection/scroll position where the edit ended up, and give the new
#    sheet its own remembered selection.
# ---------------------------------------------------------------------
$notFound.Activate()
$notFound.Range("A4").Select()

$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 231
$aw.ScrollColumn = 1
$ws.Range("A242").Select()
